# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.984.31"
$ws.Range("E2").Value = "  -2.43%  "

$ws.Range("D3").Value = "1.663.41"
$ws.Range("E3").Value = "  -2.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5088"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.006"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2631"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06384"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07412"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("D12").Value = "1.666.37"
$ws.Range("E12").Value = "  -1.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.495"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5804"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008468"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.30%  "

$ws.Range("D17").Value = "26.054.29"
$ws.Range("E17").Value = "  -2.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.902"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.90%  "

$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.193"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.007"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.570"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.61%  "

$ws.Range("E26").Value = "  +2.99%  "

$ws.Range("E27").Value = "  -1.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06671"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +16.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.302"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.54%  "

$ws.Range("E30").Value = "  -1.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.520"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.494"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.623"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.017"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6052"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.368"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.680"
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.210"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01608"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("D40").Value = "1.075.52"
$ws.Range("E40").Value = "  -2.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8581"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.80%  "

$ws.Range("E42").Value = "  +0.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.50%  "

$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.00000000116"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.20%  "

$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.812.49"
$ws.Range("E45").Value = "  -2.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.005"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.996"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05208"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4289"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.941"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.77%  "
